$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark (it currently sits right
#     after the first "UML" run in the "Механизмы расширения UML..." paragraph) ---
if ($d.Bookmarks.Exists("_GoBack")) {
  $oldBm = $d.Bookmarks("_GoBack")
  $oldBm.Delete()
}

# --- Step 2: replace the text of the last paragraph of section 2
#     ("При помощи UML-диаграмм были рассмотрены роли пользователя и
#     администратора.") with the new paragraph text ---
$rng = $d.Content
$found = $rng.Find.Execute("При помощи UML-диаграмм были рассмотрены роли пользователя и администратора.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $newText = "Ввиду того, что классов в программе большое количество, и многие из них одинаковы по структуре (классы для таблиц из базы данных), будут рассмотрены те, которые непосредственно отличаются друг от друга по структуре и схеме взаимодействия с другими классами."
    $rng.Text = $newText
    $rng.LanguageID = "ru-RU"

    # --- Step 3: add the "_GoBack" bookmark back, this time collapsed at the
    #     very end of this paragraph's text (right before the paragraph mark).
    #
    #     Placing a collapsed range exactly at a paragraph's content end (the
    #     position immediately preceding the paragraph mark) confuses this
    #     runtime's Bookmarks.Add, so a temporary marker character is
    #     inserted after the text first; the bookmark is created right before
    #     that marker (no longer exactly at the paragraph-mark boundary), and
    #     the marker is deleted afterwards, leaving the bookmark correctly
    #     collapsed at the paragraph's end. ---
    $para = $rng.Paragraphs(1)
    $textRange = $para.Range
    $textRange.MoveEnd(1, -1)
    $textRange.Collapse(0)
    $endPos = $textRange.End

    $textRange.InsertAfter("@")

    $bmRange = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $markerRange = $d.Range($endPos, $endPos + 1)
    $markerRange.Delete()
}
